$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Tabelle1" to "Table1"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Table1"

# Add a new worksheet right after Table1 and name it "Table2"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Table2"

# Populate Table2 with the same data as Table1 (A1:B3)
$ws2.Range("A1").Value = "product_ID"
$ws2.Range("B1").Value = "username"
$ws2.Range("A2").Value = 253
$ws2.Range("B2").Value = "testUser"
$ws2.Range("A3").Value = 254
$ws2.Range("B3").Value = "testUser2"

# Select A1:B3 on the new sheet
$ws2.Range("A1:B3").Select()

# Go back to Table1 and move the selection to B19
$ws1.Activate()
$ws1.Range("B19").Select()
